# Convert the "RM..." price-text cells in columns E and F (rows 131-176)
# into plain numeric values. The cells already carry a custom number
# format ("RM"#,##0.00) so the "RM..." prefix continues to display once
# the underlying value is a real number instead of literal text.
#
# Removing these text values also drops the now-unreferenced shared
# strings automatically on save, which re-indexes every other shared
# string reference throughout the sheet (e.g. column H's "Clear On hand
# Stock" pointer) without any extra work here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prices = @{
    131 = 14.4
    132 = 15.5
    133 = 16.9
    134 = 18.9
    135 = 9.9
    136 = 18.5
    137 = 11.9
    138 = 18
    139 = 18.5
    140 = 4.5
    141 = 5.7
    142 = 3.5
    147 = 3.3
    148 = 11.5
    149 = 8.1
    150 = 2.6
    151 = 9.1
    152 = 16.5
    153 = 7.9
    170 = 30
    171 = 36
    172 = 35
    173 = 35
    174 = 15
    175 = 40
    176 = 50
}

foreach ($row in $prices.Keys) {
    $val = $prices[$row]
    $ws.Range("E$row").Value = $val
    $ws.Range("F$row").Value = $val
}

# Restore the view state left by the author when the workbook was saved
# (scrolled roughly to the middle of the sheet, with E113 selected).
[void]$ws.Range("A86").Select()
[void]$ws.Range("E113").Select()
